# Fixed variables and query errors in Bread from TC01 to TC30
#
# The "CasesTab" query (cell B2 on the "startup" sheet) incorrectly pulled in a
# cohort/cohort_description match+column that doesn't belong to the Cases query.
# This removes the extra OPTIONAL MATCH (co:cohort)... plumbing's trailing output
# column (`Cohort`) from the end of that query text, restoring the cell to the
# correct query (File/Association/etc. columns for other tabs are unaffected).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("startup")

$casesQuery = @'
MATCH (s:study)<-[*]-(c:case)<--(demo:demographic)
WHERE demo.breed IN ['Chinese Shar-Pei']
MATCH (c)<--(diag:diagnosis)
OPTIONAL MATCH (samp:sample)-->(c)
OPTIONAL MATCH (co:cohort)<-[*]-(c)
WITH DISTINCT c, s, demo, diag, co
RETURN  coalesce(c.case_id, '') AS `Case ID` ,
        coalesce(s.clinical_study_designation, '') AS `Study Code` ,
        coalesce(s.clinical_study_type, '') AS  `Study Type`,
        coalesce(demo.breed, '') AS Breed ,
        coalesce(diag.disease_term, '') AS Diagnosis ,
        coalesce(diag.stage_of_disease, '') AS `Stage of Disease` ,
        coalesce(demo.patient_age_at_enrollment, '') AS Age ,
        coalesce(demo.sex, '') AS Sex ,
        coalesce(demo.neutered_indicator, '') AS `Neutered Status`,
        coalesce(demo.weight, '') AS `Weight (kg)`,
        coalesce(diag.best_response, '') AS `Response to Treatment`
'@

$ws.Range("B2").Value = $casesQuery

# Row heights settle to Excel's recalculated wrap-text autofit values once the
# query text shrinks.
$ws.Rows.Item(2).RowHeight = 244.8
$ws.Rows.Item(3).RowHeight = 230.4
$ws.Rows.Item(4).RowHeight = 216

# Restore the view to the top of the sheet with B2 (CasesTab row) selected,
# instead of being scrolled down to / focused on B4.
$ws.Range("B2").Select()
